$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the numeric-looking "price" and "volume" columns (D and E) keep their
# original plain-text representation instead of being auto-converted to numbers
# or percentages by Excel when the new values are assigned.
$deCells = @("D2", "E2", "D3", "E3", "D4", "E4", "D5", "E5", "D6", "E6", "D7", "E7", "D8", "E8", "D9", "E9", "D10", "E10", "D11", "E11", "D12", "E12", "D13", "E13", "D14", "E14", "D15", "E15", "D16", "E16", "D17", "E17", "E18", "D19", "E19", "D20", "E20", "D21", "E21", "D22", "E22", "D23", "D24", "E24", "E25", "D38", "E38", "D39", "E39", "D40", "E40", "D41", "E41", "D42", "E42", "E43", "D44", "E44", "D45", "E45", "D46", "E46", "E47", "D48", "E48", "D49", "E49", "E50", "E51")
foreach ($ref in $deCells) {
    $ws.Range($ref).NumberFormat = "@"
}

$ws.Range("D2").Value = "307.89"
$ws.Range("E2").Value = "0.52%"
$ws.Range("D3").Value = "41.11"
$ws.Range("E3").Value = "3.13%"
$ws.Range("D4").Value = "5.123"
$ws.Range("E4").Value = "0.82%"
$ws.Range("D5").Value = "0.07613"
$ws.Range("E5").Value = "-1.04%"
$ws.Range("D6").Value = "1.619"
$ws.Range("E6").Value = "0.25%"
$ws.Range("B7").Value = "BTSEToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D7").Value = "2.448"
$ws.Range("E7").Value = "0.49%"
$ws.Range("B8").Value = "MXToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D8").Value = "0.8990"
$ws.Range("E8").Value = "2.31%"
$ws.Range("D9").Value = "0.1087"
$ws.Range("E9").Value = "12.64%"
$ws.Range("D10").Value = "0.1764"
$ws.Range("E10").Value = "2.37%"
$ws.Range("D11").Value = "0.09158"
$ws.Range("E11").Value = "2.82%"
$ws.Range("D12").Value = "0.04204"
$ws.Range("E12").Value = "-4.57%"
$ws.Range("D13").Value = "0.1051"
$ws.Range("E13").Value = "-0.49%"
$ws.Range("D14").Value = "0.001249"
$ws.Range("E14").Value = "-0.31%"
$ws.Range("D15").Value = "0.005844"
$ws.Range("E15").Value = "-0.84%"
$ws.Range("D16").Value = "3.352"
$ws.Range("E16").Value = "-0.14%"
$ws.Range("D17").Value = "4.254"
$ws.Range("E17").Value = "-0.43%"
$ws.Range("E18").Value = "-0.19%"
$ws.Range("D19").Value = "6.561"
$ws.Range("E19").Value = "-6.97%"
$ws.Range("D20").Value = "0.1364"
$ws.Range("E20").Value = "1.92%"
$ws.Range("D21").Value = "0.2682"
$ws.Range("E21").Value = "-16.71%"
$ws.Range("D22").Value = "0.04076"
$ws.Range("E22").Value = "-2.86%"
$ws.Range("D23").Value = "0.001225"
$ws.Range("D24").Value = "0.004095"
$ws.Range("E24").Value = "0.77%"
$ws.Range("E25").Value = "6.60%"
$ws.Range("D38").Value = "0.02376"
$ws.Range("E38").Value = "1.52%"
$ws.Range("D39").Value = "0.05178"
$ws.Range("E39").Value = "0.62%"
$ws.Range("D40").Value = "0.007772"
$ws.Range("E40").Value = "-2.24%"
$ws.Range("D41").Value = "0.1298"
$ws.Range("E41").Value = "-1.85%"
$ws.Range("D42").Value = "0.006770"
$ws.Range("E42").Value = "6.12%"
$ws.Range("E43").Value = "0.09%"
$ws.Range("D44").Value = "0.008546"
$ws.Range("E44").Value = "-0.56%"
$ws.Range("D45").Value = "0.3071"
$ws.Range("E45").Value = "0.88%"
$ws.Range("D46").Value = "0.00007033"
$ws.Range("E46").Value = "7.82%"
$ws.Range("E47").Value = "0.04%"
$ws.Range("D48").Value = "0.03414"
$ws.Range("E48").Value = "912.67%"
$ws.Range("D49").Value = "0.004203"
$ws.Range("E49").Value = "-39.97%"
$ws.Range("E50").Value = "0.04%"
$ws.Range("E51").Value = "0.04%"

Write-Output "Applied 73 cell updates"
